$wb = $excel.ActiveWorkbook

$wsExt = $wb.Worksheets.Item("ExtensionSchemes_exttest1")
$wsMembers = $wb.Worksheets.Item("Extensions_test")
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")

# Sheet3 ("ExtensionSchemes_exttest1" -> later "Extensions_exttest1"): update H1/H2 reference labels
$wsExt.Range("H1").Value = "MEMBERSSHEET"
$wsExt.Range("H2").Value = "Members_test"

# Sheet4 ("Extensions_test" -> later "Members_test"): add "code:" prefix to D3:D11
for ($r = 3; $r -le 11; $r++) {
    $cell = $wsMembers.Cells.Item($r, 4)
    $cell.Value = "code:" + $cell.Value2
}

# Sheet1 (CodeSchemes): update N1/N2 reference labels
$wsCodeSchemes.Range("N1").Value = "EXTENSIONSSHEET"
$wsCodeSchemes.Range("N2").Value = "Extensions_exttest1"

# Rename sheets
$wsExt.Name = "Extensions_exttest1"
$wsMembers.Name = "Members_test"
